$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Update 'Price' column (D) ---
# Force text format first so numeric-looking strings (e.g. '551.60')
# are not auto-converted to numbers by Excel, matching the original
# inline-string cell type/content.
$dRange = $ws.Range("D2:D48")
$dRange.NumberFormat = "@"

$ws.Range("D2").Value = "60.037.36"
$ws.Range("D3").Value = "2.416.68"
$ws.Range("D5").Value = "551.60"
$ws.Range("D6").Value = "137.04"
$ws.Range("D8").Value = "0.597"
$ws.Range("D13").Value = "25.32"
$ws.Range("D14").Value = "2.849.19"
$ws.Range("D15").Value = "60.001.91"
$ws.Range("D17").Value = "2.423.20"
$ws.Range("D18").Value = "11.29"
$ws.Range("D20").Value = "327.95"
$ws.Range("D21").Value = "6.66"
$ws.Range("D23").Value = "65.86"
$ws.Range("D24").Value = "0.177"
$ws.Range("D25").Value = "8.59"
$ws.Range("D30").Value = "168.82"
$ws.Range("D31").Value = "6.04"
$ws.Range("D32").Value = "18.56"
$ws.Range("D33").Value = "1.03"
$ws.Range("D35").Value = "1.31"
$ws.Range("D39").Value = "321.79"
$ws.Range("D42").Value = "139.85"
$ws.Range("D43").Value = "0.0973"
$ws.Range("D45").Value = "0.0514"
$ws.Range("D47").Value = "0.0224"
$ws.Range("D48").Value = "0.386"

# Restore the default (Normal) cell style so no stray number-format
# style index is left attached to these cells.
$dRange.Style = "Normal"

# --- Update 'Volume(1h)' column (E) ---
$ws.Range("E2").Value = "  -0.42%  "
$ws.Range("E3").Value = "  -1.35%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("E5").Value = "  -0.90%  "
$ws.Range("E6").Value = "  -1.53%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("E8").Value = "  +4.20%  "
$ws.Range("E9").Value = "  -1.92%  "
$ws.Range("E10").Value = "  -2.80%  "
$ws.Range("E11").Value = "  -0.99%  "
$ws.Range("E12").Value = "  -2.50%  "
$ws.Range("E13").Value = "  +1.17%  "
$ws.Range("E14").Value = "  -1.24%  "
$ws.Range("E15").Value = "  -0.30%  "
$ws.Range("E16").Value = "  -2.39%  "
$ws.Range("E17").Value = "  -0.66%  "
$ws.Range("E18").Value = "  -1.99%  "
$ws.Range("E19").Value = "  -0.64%  "
$ws.Range("E20").Value = "  -2.38%  "
$ws.Range("E21").Value = "  -3.78%  "
$ws.Range("E22").Value = "  +0.03%  "
$ws.Range("E23").Value = "  +1.55%  "
$ws.Range("E24").Value = "  +3.98%  "
$ws.Range("E25").Value = "  -0.54%  "
$ws.Range("E26").Value = "  +0.18%  "
$ws.Range("E27").Value = "  -0.62%  "
$ws.Range("E28").Value = "  -3.12%  "
$ws.Range("E29").Value = "  -2.30%  "
$ws.Range("E30").Value = "  -1.52%  "
$ws.Range("E31").Value = "  -4.47%  "
$ws.Range("E32").Value = "  -1.70%  "
$ws.Range("E33").Value = "  -0.22%  "
$ws.Range("E34").Value = "  +0.00%  "
$ws.Range("E35").Value = "  -1.14%  "
$ws.Range("E36").Value = "  +0.02%  "
$ws.Range("E37").Value = "  -3.04%  "
$ws.Range("E38").Value = "  -2.59%  "
$ws.Range("E39").Value = "  +1.52%  "
$ws.Range("E40").Value = "  -3.74%  "
$ws.Range("E41").Value = "  -2.29%  "
$ws.Range("E42").Value = "  -3.16%  "
$ws.Range("E43").Value = "  +0.95%  "
$ws.Range("E44").Value = "  -0.95%  "
$ws.Range("E45").Value = "  -2.45%  "
$ws.Range("E46").Value = "  +0.04%  "
$ws.Range("E47").Value = "  -1.58%  "
$ws.Range("E48").Value = "  -6.12%  "
$ws.Range("E49").Value = "  +0.03%  "
$ws.Range("E50").Value = "  -4.20%  "
$ws.Range("E51").Value = "  -1.02%  "
